$d = $word.ActiveDocument

# Locate the paragraph describing the "new window per button" rationale -
# it currently contains the old (un-split) wording and still holds the
# stray _GoBack bookmark that cuts "Un joueur de " away from
# "jeu de rôle ...".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Nous avons fait le choix")) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the target paragraph ('Nous avons fait le choix ...')."
}

# Replace that single paragraph with two paragraphs:
#  1) the same rationale text, now written without the bookmark splitting
#     "Un joueur de " / "jeu de rôle ..." (still keeping the trailing
#     "d'avoir devant lui ..." run separate, as it already was);
#  2) a brand-new paragraph explaining the double-click choice for
#     displaying character details, with the (moved) _GoBack bookmark now
#     sitting inside "voul|ait" where the author was last typing.
$xmlFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Nous avons fait le choix d’ouvrir une nouvelle fenêtre à chaque fois que l’ou appuie sur un bouton. Ce choix peut paraitre très peu instinctif mais tout au contraire, ce choix est réfléchi. Un joueur de jeu de rôle a souvent besoin de voir plusieurs compétences en même temps. De plus il peu déplacé cette fenêtre ou il le souhait sur son écran. De plus, le meneur de jeu (celui qui dirige la partie) a besoin </w:t></w:r><w:r><w:t xml:space="preserve">d’avoir devant lui plusieurs profils en même temps. </w:t></w:r></w:p><w:p><w:r><w:t>Nous avons décidé d’utiliser des doubles cliques pour pouvoir afficher le détail des personnages car l’utilisateur peut sans le faire exprès cliquer sur un profil alors qu’il voul</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ait pas voir son détail.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xmlFragment)

Write-Host "Paragraph split done. Total paragraphs now:" $d.Paragraphs.Count
